# TMTC0032668 - Updated Test Data for LV Activities (16 Sep 2024)
# The "Contact" sheet's sample row is refreshed to use a dedicated
# Activity-test contact/company instead of the generic standard one,
# and the Contact sheet becomes the active/selected tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Contact")

# Update the sample contact row used by the Activity external-attendee tests.
$ws.Range("A2").Value = "Activity Test External Contact"
$ws.Range("B2").Value = "ActivityCompany"

# Make "Contact" the active sheet/tab with A2:B2 selected (matches the
# recorded workbook view state after the edit).
$ws.Activate()
$ws.Range("A2:B2").Select()
